$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions (G1, H1, J1, K1) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2 additions ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Re-enter the D3:D9 segment formula as a single fill so Excel
#     regroups it into one shared formula (matches author re-confirming
#     the existing formula over the column) ---
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- New G column (Area per segment), rows 3 standalone then 4:15 as a
#     single fill/shared-formula group ---
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Selection as per diff ---
$ws.Range("G1:K15").Select()
